$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the separate card-detail rows (A2:A7) into a single Python tuple-repr string in A2
$combined = "('Feral Throwback', ['{4}{G}{G}', 'Creature " + [char]0x2014 + " Beast', 'Amplify 2 (As this creature enters the battlefield, put two +1/+1 counters on it for each Beast card you reveal in your hand.)', 'Provoke (Whenever this creature attacks, you may have target creature defending player controls untap and block it if able.)', '3/3'])"

# Remove the now-unneeded rows 3-7 (shifts nothing below them, they're the last rows)
$ws.Range("A3:A7").EntireRow.Delete() | Out-Null

# Set the new combined value in A2
$ws.Range("A2").Value = $combined
